$d = $word.ActiveDocument

# 1) "lysines" -> "lysine"  (stays inside the existing spellStart/spellEnd pair)
$d.Content.Find.Execute("lysines", $true, $false, $false, $false, $false, $true, 1, $false, "lysine", 2) | Out-Null

# 2) "Class I, and II HDACs including class II HDAC " -> "Class I, II, III HDACs including "
#    (leaves the following "Sirtuins" run - and its spellStart/spellEnd wrapper - untouched)
$d.Content.Find.Execute("Class I, and II HDACs including class II HDAC ", $true, $false, $false, $false, $false, $true, 1, $false, "Class I, II, III HDACs including ", 2) | Out-Null

# 3) "iPSC, H3" -> "iPSC; H3"
$d.Content.Find.Execute("iPSC, H3", $true, $false, $false, $false, $false, $true, 1, $false, "iPSC; H3", 2) | Out-Null

# 4) "...chromatin regions. They include" -> "...chromatin regions, and include"
$d.Content.Find.Execute("chromatin regions. They include", $true, $false, $false, $false, $false, $true, 1, $false, "chromatin regions, and include", 2) | Out-Null

# 5) "LOCK distribution and abundance changes upon differentiation." ->
#    "LOCK distribution and abundance; change upon differentiation."
$d.Content.Find.Execute("abundance changes upon differentiation.", $true, $false, $false, $false, $false, $true, 1, $false, "abundance; change upon differentiation.", 2) | Out-Null

# 6) "...early to late replicating domains." -> "...early to late replicating domain."
$d.Content.Find.Execute("late replicating domains.", $true, $false, $false, $false, $false, $true, 1, $false, "late replicating domain.", 2) | Out-Null

# 7) Remove the trailing empty paragraph that sits right before the final empty
#    paragraph preceding the section properties. Both trailing paragraphs are
#    empty, so target the second-to-last paragraph explicitly (leaving the
#    very last paragraph of the document untouched).
$count = $d.Paragraphs.Count
$secondToLast = $d.Paragraphs.Item($count - 1)
if ($secondToLast.Range.Text -eq "`r") {
    $secondToLast.Range.Delete()
}
